$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove column M entirely (old M data is dropped, old N shifts left into M)
$ws.Columns("M").Delete()

# Update selection to reflect the new active cell after the edit
$ws.Range("M1").Select()
